# updateMapTime for TimeKeeper, volume managing
# Adds a new "numberOfWorkers" rule/value row to the rules sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row: RULE = "numberOfWorkers", VALUE = 2
$ws.Range("A6").Value = "numberOfWorkers"
$ws.Range("B6").Value = 2

# Match the integer number formatting used by the other VALUE cells (B2:B5)
$ws.Range("B6").NumberFormat = $ws.Range("B5").NumberFormat

# Leave the active selection where the user ended up after entering the row
$ws.Range("B7").Select()
